# Apply daily COVID data refresh to the "Pais" sheet.
# This mirrors the upstream dataset update: most countries keep their rank
# but get refreshed totals, while a handful of countries with close totals
# swap rank/row order once the new totals are applied.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados..." timestamp footer (A1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 19:35"

# Row 4: Estados Unidos -> Estados Unidos
$ws.Cells.Item(4, 1).Value = "Estados Unidos"
$ws.Cells.Item(4, 2).Value = 2278872
$ws.Cells.Item(4, 3).Value = 15221
$ws.Cells.Item(4, 4).Value = 932602
$ws.Cells.Item(4, 5).Value = 1225247
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 335
$ws.Cells.Item(4, 8).Value = 121023

# Row 5: Brasil -> Brasil
$ws.Cells.Item(5, 1).Value = "Brasil"
$ws.Cells.Item(5, 2).Value = 1009699
$ws.Cells.Item(5, 3).Value = 26340
$ws.Cells.Item(5, 4).Value = 520360
$ws.Cells.Item(5, 5).Value = 440912
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(5, 7).Value = 558
$ws.Cells.Item(5, 8).Value = 48427

# Row 7: India -> India
$ws.Cells.Item(7, 1).Value = "India"
$ws.Cells.Item(7, 2).Value = 392536
$ws.Cells.Item(7, 3).Value = 11445
$ws.Cells.Item(7, 4).Value = 211944
$ws.Cells.Item(7, 5).Value = 167688
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(7, 7).Value = 300
$ws.Cells.Item(7, 8).Value = 12904

# Row 14: Alemania -> Alemania
$ws.Cells.Item(14, 1).Value = "Alemania"
$ws.Cells.Item(14, 2).Value = 190444
$ws.Cells.Item(14, 3).Value = 318
$ws.Cells.Item(14, 4).Value = 174400
$ws.Cells.Item(14, 5).Value = 7092
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(14, 7).Value = 6
$ws.Cells.Item(14, 8).Value = 8952

# Row 15: Turquia -> Turquia
$ws.Cells.Item(15, 1).Value = "Turquia"
$ws.Cells.Item(15, 2).Value = 185245
$ws.Cells.Item(15, 3).Value = 1214
$ws.Cells.Item(15, 4).Value = 157516
$ws.Cells.Item(15, 5).Value = 22824
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 23
$ws.Cells.Item(15, 8).Value = 4905

# Row 21: Canada -> Canada
$ws.Cells.Item(21, 1).Value = "Canada"
$ws.Cells.Item(21, 2).Value = 100565
$ws.Cells.Item(21, 3).Value = 345
$ws.Cells.Item(21, 4).Value = 62961
$ws.Cells.Item(21, 5).Value = 29258
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 46
$ws.Cells.Item(21, 8).Value = 8346

# Row 44: Irak -> Irak
$ws.Cells.Item(44, 1).Value = "Irak"
$ws.Cells.Item(44, 2).Value = 27352
$ws.Cells.Item(44, 3).Value = 1635
$ws.Cells.Item(44, 4).Value = 12205
$ws.Cells.Item(44, 5).Value = 14222
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 69
$ws.Cells.Item(44, 8).Value = 925

# Row 45: Irlanda -> Irlanda
$ws.Cells.Item(45, 1).Value = "Irlanda"
$ws.Cells.Item(45, 2).Value = 25368
$ws.Cells.Item(45, 3).Value = 13
$ws.Cells.Item(45, 4).Value = 22698
$ws.Cells.Item(45, 5).Value = 956
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(45, 7).Value = 0
$ws.Cells.Item(45, 8).Value = 1714

# Row 51: Israel -> Israel
$ws.Cells.Item(51, 1).Value = "Israel"
$ws.Cells.Item(51, 2).Value = 20339
$ws.Cells.Item(51, 3).Value = 303
$ws.Cells.Item(51, 4).Value = 15586
$ws.Cells.Item(51, 5).Value = 4449
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 304

# Row 68: Marruecos -> Marruecos
$ws.Cells.Item(68, 1).Value = "Marruecos"
$ws.Cells.Item(68, 2).Value = 9613
$ws.Cells.Item(68, 3).Value = 539
$ws.Cells.Item(68, 4).Value = 8117
$ws.Cells.Item(68, 5).Value = 1283
$ws.Cells.Item(68, 6).Value = 0
$ws.Cells.Item(68, 7).Value = 0
$ws.Cells.Item(68, 8).Value = 213

# Row 71: Nepal -> Sudan
$ws.Cells.Item(71, 1).Value = "Sudan"
$ws.Cells.Item(71, 2).Value = 8316
$ws.Cells.Item(71, 3).Value = 296
$ws.Cells.Item(71, 4).Value = 3086
$ws.Cells.Item(71, 5).Value = 4724
$ws.Cells.Item(71, 6).Value = 0
$ws.Cells.Item(71, 7).Value = 19
$ws.Cells.Item(71, 8).Value = 506

# Row 72: Sudan -> Nepal
$ws.Cells.Item(72, 1).Value = "Nepal"
$ws.Cells.Item(72, 2).Value = 8274
$ws.Cells.Item(72, 3).Value = 426
$ws.Cells.Item(72, 4).Value = 1402
$ws.Cells.Item(72, 5).Value = 6850
$ws.Cells.Item(72, 6).Value = 0
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 22

# Row 83: Republica de Yibuti -> Republica de Yibuti
$ws.Cells.Item(83, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(83, 2).Value = 4565
$ws.Cells.Item(83, 3).Value = 8
$ws.Cells.Item(83, 4).Value = 3565
$ws.Cells.Item(83, 5).Value = 955
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 2
$ws.Cells.Item(83, 8).Value = 45

# Row 102: Maldivas -> Maldivas
$ws.Cells.Item(102, 1).Value = "Maldivas"
$ws.Cells.Item(102, 2).Value = 2150
$ws.Cells.Item(102, 3).Value = 13
$ws.Cells.Item(102, 4).Value = 1769
$ws.Cells.Item(102, 5).Value = 373
$ws.Cells.Item(102, 6).Value = 0
$ws.Cells.Item(102, 7).Value = 0
$ws.Cells.Item(102, 8).Value = 8

# Row 136: Cabo Verde -> Cabo Verde
$ws.Cells.Item(136, 1).Value = "Cabo Verde"
$ws.Cells.Item(136, 2).Value = 849
$ws.Cells.Item(136, 3).Value = 26
$ws.Cells.Item(136, 4).Value = 377
$ws.Cells.Item(136, 5).Value = 464
$ws.Cells.Item(136, 6).Value = 0
$ws.Cells.Item(136, 7).Value = 1
$ws.Cells.Item(136, 8).Value = 8

# Row 142: Estado de Palestina -> Estado de Palestina
$ws.Cells.Item(142, 1).Value = "Estado de Palestina"
$ws.Cells.Item(142, 2).Value = 665
$ws.Cells.Item(142, 3).Value = 65
$ws.Cells.Item(142, 4).Value = 437
$ws.Cells.Item(142, 5).Value = 225
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 3

# Row 144: Ruanda -> Benin
$ws.Cells.Item(144, 1).Value = "Benin"
$ws.Cells.Item(144, 2).Value = 650
$ws.Cells.Item(144, 3).Value = 53
$ws.Cells.Item(144, 4).Value = 247
$ws.Cells.Item(144, 5).Value = 392
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 11

# Row 145: Jamaica -> Ruanda
$ws.Cells.Item(145, 1).Value = "Ruanda"
$ws.Cells.Item(145, 2).Value = 646
$ws.Cells.Item(145, 3).Value = 0
$ws.Cells.Item(145, 4).Value = 350
$ws.Cells.Item(145, 5).Value = 294
$ws.Cells.Item(145, 6).Value = 0
$ws.Cells.Item(145, 7).Value = 0
$ws.Cells.Item(145, 8).Value = 2

# Row 146: Benin -> Jamaica
$ws.Cells.Item(146, 1).Value = "Jamaica"
$ws.Cells.Item(146, 2).Value = 638
$ws.Cells.Item(146, 3).Value = 12
$ws.Cells.Item(146, 4).Value = 458
$ws.Cells.Item(146, 5).Value = 170
$ws.Cells.Item(146, 6).Value = 0
$ws.Cells.Item(146, 7).Value = 0
$ws.Cells.Item(146, 8).Value = 10

# Row 202: Fiyi -> Dominica
$ws.Cells.Item(202, 1).Value = "Dominica"
$ws.Cells.Item(202, 2).Value = 18
$ws.Cells.Item(202, 3).Value = 0
$ws.Cells.Item(202, 4).Value = 18
$ws.Cells.Item(202, 5).Value = 0
$ws.Cells.Item(202, 6).Value = 0
$ws.Cells.Item(202, 7).Value = 0
$ws.Cells.Item(202, 8).Value = 0

# Row 203: Dominica -> Fiyi
$ws.Cells.Item(203, 1).Value = "Fiyi"
$ws.Cells.Item(203, 2).Value = 18
$ws.Cells.Item(203, 3).Value = 0
$ws.Cells.Item(203, 4).Value = 18
$ws.Cells.Item(203, 5).Value = 0
$ws.Cells.Item(203, 6).Value = 0
$ws.Cells.Item(203, 7).Value = 0
$ws.Cells.Item(203, 8).Value = 0

# Row 206: Groenlandia -> Islas Malvinas
$ws.Cells.Item(206, 1).Value = "Islas Malvinas"
$ws.Cells.Item(206, 2).Value = 13
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 13
$ws.Cells.Item(206, 5).Value = 0
$ws.Cells.Item(206, 6).Value = 0
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 0

# Row 207: Islas Malvinas -> Groenlandia
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(207, 2).Value = 13
$ws.Cells.Item(207, 3).Value = 0
$ws.Cells.Item(207, 4).Value = 13
$ws.Cells.Item(207, 5).Value = 0
$ws.Cells.Item(207, 6).Value = 0
$ws.Cells.Item(207, 7).Value = 0
$ws.Cells.Item(207, 8).Value = 0

# Row 210: Montserrat -> Seychelles
$ws.Cells.Item(210, 1).Value = "Seychelles"
$ws.Cells.Item(210, 2).Value = 11
$ws.Cells.Item(210, 3).Value = 0
$ws.Cells.Item(210, 4).Value = 11
$ws.Cells.Item(210, 5).Value = 0
$ws.Cells.Item(210, 6).Value = 0
$ws.Cells.Item(210, 7).Value = 0
$ws.Cells.Item(210, 8).Value = 0

# Row 211: Seychelles -> Montserrat
$ws.Cells.Item(211, 1).Value = "Montserrat"
$ws.Cells.Item(211, 2).Value = 11
$ws.Cells.Item(211, 3).Value = 0
$ws.Cells.Item(211, 4).Value = 10
$ws.Cells.Item(211, 5).Value = 0
$ws.Cells.Item(211, 6).Value = 0
$ws.Cells.Item(211, 7).Value = 0
$ws.Cells.Item(211, 8).Value = 1

# Row 213: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(213, 2).Value = 8
$ws.Cells.Item(213, 3).Value = 0
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 5).Value = 0
$ws.Cells.Item(213, 6).Value = 0
$ws.Cells.Item(213, 7).Value = 0
$ws.Cells.Item(213, 8).Value = 0

# Row 214: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"
$ws.Cells.Item(214, 2).Value = 8
$ws.Cells.Item(214, 3).Value = 0
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 5).Value = 0
$ws.Cells.Item(214, 6).Value = 0
$ws.Cells.Item(214, 7).Value = 0
$ws.Cells.Item(214, 8).Value = 1

Write-Host "Applied country data refresh for 19 Jun 2020 19:35"